$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8535362482070923
$ws.Range("B1").Value = 2.635683059692383
$ws.Range("C1").Value = 1.128680467605591
$ws.Range("D1").Value = 1.188580393791199
$ws.Range("E1").Value = 1.346000075340271
